$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws3 = $wb.Worksheets.Item("Final")

# --- Sheet "2o Parcial": rows 2-17, columns E:K (Aprobados..por_blancos) ---
$ws2.Range("E2").Value = 35
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 100
$ws2.Range("H2").Value = 0
$ws2.Range("I2").Value = 8.9
$ws2.Range("J2").Value = 0
$ws2.Range("K2").Value = 0
$ws2.Range("E3").Value = 39
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 100
$ws2.Range("H3").Value = 0
$ws2.Range("I3").Value = 9.199999999999999
$ws2.Range("J3").Value = 0
$ws2.Range("K3").Value = 0
$ws2.Range("E4").Value = 28
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 100
$ws2.Range("H4").Value = 0
$ws2.Range("I4").Value = 8.300000000000001
$ws2.Range("J4").Value = 0
$ws2.Range("K4").Value = 0
$ws2.Range("E5").Value = 25
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 100
$ws2.Range("H5").Value = 0
$ws2.Range("I5").Value = 9.199999999999999
$ws2.Range("J5").Value = 0
$ws2.Range("K5").Value = 0
$ws2.Range("E6").Value = 127
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 100
$ws2.Range("H6").Value = 0
$ws2.Range("I6").Value = 8.9
$ws2.Range("J6").Value = 0
$ws2.Range("K6").Value = 0
$ws2.Range("E7").Value = 13
$ws2.Range("F7").Value = 7
$ws2.Range("G7").Value = 65
$ws2.Range("H7").Value = 35
$ws2.Range("I7").Value = 6
$ws2.Range("J7").Value = 0
$ws2.Range("K7").Value = 0
$ws2.Range("E8").Value = 13
$ws2.Range("F8").Value = 7
$ws2.Range("G8").Value = 65
$ws2.Range("H8").Value = 35
$ws2.Range("I8").Value = 6
$ws2.Range("J8").Value = 0
$ws2.Range("K8").Value = 0
$ws2.Range("E9").Value = 27
$ws2.Range("F9").Value = 1
$ws2.Range("G9").Value = 96.40000000000001
$ws2.Range("H9").Value = 3.6
$ws2.Range("I9").Value = 8.1
$ws2.Range("J9").Value = 0
$ws2.Range("K9").Value = 0
$ws2.Range("E10").Value = 25
$ws2.Range("F10").Value = 6
$ws2.Range("G10").Value = 80.59999999999999
$ws2.Range("H10").Value = 19.4
$ws2.Range("I10").Value = 7.3
$ws2.Range("J10").Value = 0
$ws2.Range("K10").Value = 0
$ws2.Range("E11").Value = 52
$ws2.Range("F11").Value = 7
$ws2.Range("G11").Value = 88.09999999999999
$ws2.Range("H11").Value = 11.9
$ws2.Range("I11").Value = 7.7
$ws2.Range("J11").Value = 0
$ws2.Range("K11").Value = 0
$ws2.Range("E12").Value = 40
$ws2.Range("F12").Value = 1
$ws2.Range("G12").Value = 97.59999999999999
$ws2.Range("H12").Value = 2.4
$ws2.Range("I12").Value = 8.9
$ws2.Range("J12").Value = 0
$ws2.Range("K12").Value = 0
$ws2.Range("E13").Value = 35
$ws2.Range("F13").Value = 1
$ws2.Range("G13").Value = 97.2
$ws2.Range("H13").Value = 2.8
$ws2.Range("I13").Value = 9
$ws2.Range("J13").Value = 0
$ws2.Range("K13").Value = 0
$ws2.Range("E14").Value = 17
$ws2.Range("F14").Value = 0
$ws2.Range("G14").Value = 100
$ws2.Range("H14").Value = 0
$ws2.Range("I14").Value = 8.800000000000001
$ws2.Range("J14").Value = 0
$ws2.Range("K14").Value = 0
$ws2.Range("E15").Value = 24
$ws2.Range("F15").Value = 0
$ws2.Range("G15").Value = 100
$ws2.Range("H15").Value = 0
$ws2.Range("I15").Value = 8.9
$ws2.Range("J15").Value = 0
$ws2.Range("K15").Value = 0
$ws2.Range("E16").Value = 116
$ws2.Range("F16").Value = 2
$ws2.Range("G16").Value = 98.3
$ws2.Range("H16").Value = 1.7
$ws2.Range("I16").Value = 8.9
$ws2.Range("J16").Value = 0
$ws2.Range("K16").Value = 0
$ws2.Range("E17").Value = 308
$ws2.Range("F17").Value = 16
$ws2.Range("G17").Value = 95.09999999999999
$ws2.Range("H17").Value = 4.9
$ws2.Range("I17").Value = 8.4
$ws2.Range("J17").Value = 0
$ws2.Range("K17").Value = 0

# --- Sheet "Final": selective cell updates ---
$ws3.Range("I2").Value = 9.300000000000001
$ws3.Range("E3").Value = 39
$ws3.Range("F3").Value = 0
$ws3.Range("G3").Value = 100
$ws3.Range("H3").Value = 0
$ws3.Range("I3").Value = 9
$ws3.Range("I4").Value = 8.699999999999999
$ws3.Range("E6").Value = 127
$ws3.Range("F6").Value = 0
$ws3.Range("G6").Value = 100
$ws3.Range("H6").Value = 0
$ws3.Range("I6").Value = 9.199999999999999
$ws3.Range("E7").Value = 13
$ws3.Range("F7").Value = 7
$ws3.Range("G7").Value = 65
$ws3.Range("H7").Value = 35
$ws3.Range("I7").Value = 6
$ws3.Range("E8").Value = 13
$ws3.Range("F8").Value = 7
$ws3.Range("G8").Value = 65
$ws3.Range("H8").Value = 35
$ws3.Range("I8").Value = 6
$ws3.Range("I10").Value = 7.8
$ws3.Range("I11").Value = 8.199999999999999
$ws3.Range("I12").Value = 8.699999999999999
$ws3.Range("I13").Value = 9.1
$ws3.Range("I14").Value = 9
$ws3.Range("E15").Value = 24
$ws3.Range("F15").Value = 0
$ws3.Range("G15").Value = 100
$ws3.Range("H15").Value = 0
$ws3.Range("I15").Value = 8.800000000000001
$ws3.Range("E16").Value = 116
$ws3.Range("F16").Value = 2
$ws3.Range("G16").Value = 98.3
$ws3.Range("H16").Value = 1.7
$ws3.Range("I16").Value = 8.9
$ws3.Range("E17").Value = 308
$ws3.Range("F17").Value = 16
$ws3.Range("G17").Value = 95.09999999999999
$ws3.Range("H17").Value = 4.9
$ws3.Range("I17").Value = 8.6
